$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.854.65'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '3.134.50'
$ws.Range("E3").Value = '  +3.22%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '3.126.35'
$ws.Range("E8").Value = '  +3.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("E10").Value = '  +13.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.158'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '3.636.92'
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("D16").Value = '64.866.48'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '524.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +10.79%  '
$ws.Range("D19").Value = '3.135.92'
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.86%  '
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.70%  '
$ws.Range("E28").Value = '  +3.29%  '
$ws.Range("E29").Value = '  +2.89%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '558.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.56%  '
$ws.Range("E35").Value = '  +2.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.63%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0438'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.67%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0811'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.79%  '
$ws.Range("D40").Value = '3.056.93'
$ws.Range("E40").Value = '  +7.02%  '
$ws.Range("E41").Value = '  +15.68%  '
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.255'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.64%  '
$ws.Range("D49").Value = '0.0₃0521'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("E51").Value = '  +4.92%  '
